$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "data dictionary" rows (2-5) explaining each header column.
# Cells are written in the exact order needed so that new shared strings are
# appended to the shared string table in the same order the original author's
# edit produced them.
$ws.Range("F2").Value = "bottom temp"
$ws.Range("A2").Value = "number"
$ws.Range("B2").Value = "year, sample #, month"
$ws.Range("B3").Value = "unique id for space & time of haul"
$ws.Range("D2").Value = "weight catch per unit effort - measure of biomass"
$ws.Range("D3").Value = "effort: area of sf sampled"
$ws.Range("D4").Value = "catch: kg"
$ws.Range("D5").Value = "kg/area"
$ws.Range("E2").Value = "surface temp"
$ws.Range("G2").Value = "m"
$ws.Range("J2").Value = "region identifier - ebs"
$ws.Range("N2").Value = "combination of long lat depth"
$ws.Range("P2").Value = "same as spp"
$ws.Range("Q2").Value = "1 in most cases"
$ws.Range("Q3").Value = "number of samples per place & time"
$ws.Range("R2").Value = "# of indiv"
$ws.Range("R3").Value = "don't read into"
$ws.Range("Q4").Value = "ignore"
$ws.Range("R4").Value = "ignore"
$ws.Range("T3").Value = "number of hauls w/in stratum year combo"
$ws.Range("S4").Value = "indicator of replicates in stratum year"
$ws.Range("S2").Value = "1 to 4"
$ws.Range("S2").NumberFormat = "d-mmm"

# C2 re-uses the already existing "species" shared string.
$ws.Range("C2").Value = "species"

# Column widths were adjusted (best-fit) after the new columns of descriptive
# text were added.
$ws.Columns("B:B").ColumnWidth = 26.75014
$ws.Columns("C:C").ColumnWidth = 5.4167
$ws.Columns("D:D").ColumnWidth = 39.08357
$ws.Columns("E:E").ColumnWidth = 10.08367
$ws.Columns("F:F").ColumnWidth = 10.25039
$ws.Columns("J:J").ColumnWidth = 16.75014
$ws.Columns("N:N").ColumnWidth = 23.4167
$ws.Columns("P:P").ColumnWidth = 9.25039
$ws.Columns("Q:Q").ColumnWidth = 28.75032
$ws.Columns("R:R").ColumnWidth = 11.4171

# Final selection, matching where the user's cursor ended up.
$ws.Range("D6").Select()
